$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.051956481319977
$ws.Range("D2").Value = 1.056948362040533
$ws.Range("E2").Value = 1.058866272110174
$ws.Range("F2").Value = 1.069642343597047
$ws.Range("I2").Value = 1.047788579277229
$ws.Range("J2").Value = 1.056981621026416
$ws.Range("K2").Value = 1.059684673658182
$ws.Range("L2").Value = 1.061597336675522
$ws.Range("M2").Value = 1.072344305255365
$ws.Range("N2").Value = 1.02277518295572
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.052893812425761
$ws.Range("D3").Value = 1.057671159107791
$ws.Range("E3").Value = 1.05968258229249
$ws.Range("F3").Value = 1.070497694383154
$ws.Range("I3").Value = 1.048031740121935
$ws.Range("J3").Value = 1.057569302264864
$ws.Range("K3").Value = 1.060221441744164
$ws.Range("L3").Value = 1.062227758498216
$ws.Range("M3").Value = 1.073015768316016
$ws.Range("N3").Value = 1.022973954822337
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.053500950112082
$ws.Range("D4").Value = 1.058139404622995
$ws.Range("E4").Value = 1.060211682419284
$ws.Range("F4").Value = 1.071052095394256
$ws.Range("I4").Value = 1.04818824511438
$ws.Range("J4").Value = 1.057949541674133
$ws.Range("K4").Value = 1.060568640864415
$ws.Range("L4").Value = 1.062635921745681
$ws.Range("M4").Value = 1.073450533347499
$ws.Range("N4").Value = 1.023102485099783
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053756338394634
$ws.Range("D5").Value = 1.058336384569341
$ws.Range("E5").Value = 1.060434328280976
$ws.Range("F5").Value = 1.071285387009342
$ws.Range("I5").Value = 1.048253838957543
$ws.Range("J5").Value = 1.058109386086177
$ws.Range("K5").Value = 1.060714572137483
$ws.Range("L5").Value = 1.062807569340572
$ws.Range("M5").Value = 1.073633375084245
$ws.Range("N5").Value = 1.023156497602481
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053799227829905
$ws.Range("D6").Value = 1.058369465920236
$ws.Range("E6").Value = 1.060471723896509
$ws.Range("F6").Value = 1.071324570634623
$ws.Range("I6").Value = 1.048264840655738
$ws.Range("J6").Value = 1.058136224150344
$ws.Range("K6").Value = 1.060739072775539
$ws.Range("L6").Value = 1.062836392960756
$ws.Range("M6").Value = 1.073664078877483
$ws.Range("N6").Value = 1.023165565257827
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.053504362044411
$ws.Range("D7").Value = 1.058142036171634
$ws.Range("E7").Value = 1.060214656592491
$ws.Range("F7").Value = 1.071055211780244
$ws.Range("I7").Value = 1.048189122372538
$ws.Range("J7").Value = 1.057951677556771
$ws.Range("K7").Value = 1.060570590929392
$ws.Range("L7").Value = 1.062638215091693
$ws.Range("M7").Value = 1.073452976228068
$ws.Range("N7").Value = 1.023103206903039
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.052273127254112
$ws.Range("D8").Value = 1.057192520410372
$ws.Range("E8").Value = 1.059141962135224
$ws.Range("F8").Value = 1.069931219455444
$ws.Range("I8").Value = 1.04787092940871
$ws.Range("J8").Value = 1.057180236106987
$ws.Range("K8").Value = 1.059866102293019
$ws.Range("L8").Value = 1.061810340612796
$ws.Range("M8").Value = 1.072571169814929
$ws.Range("N8").Value = 1.022842376868331
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.050108352019812
$ws.Range("D9").Value = 1.055523613275838
$ws.Range("E9").Value = 1.057258645614319
$ws.Range("F9").Value = 1.067957811617446
$ws.Range("I9").Value = 1.047303855602862
$ws.Range("J9").Value = 1.055820680950554
$ws.Range("K9").Value = 1.058623790765569
$ws.Range("L9").Value = 1.06035339955913
$ws.Range("M9").Value = 1.071019545937278
$ws.Range("N9").Value = 1.022382104219704
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.048668481993949
$ws.Range("D10").Value = 1.054413964680875
$ws.Range("E10").Value = 1.056007840004582
$ws.Range("F10").Value = 1.066647152768493
$ws.Range("I10").Value = 1.046921556581194
$ws.Range("J10").Value = 1.054914259495772
$ws.Range("K10").Value = 1.057795037163016
$ws.Range("L10").Value = 1.05938344367095
$ws.Range("M10").Value = 1.06998671133287
$ws.Range("N10").Value = 1.022074841488828
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.048045803686148
$ws.Range("D11").Value = 1.053934194947113
$ws.Range("E11").Value = 1.055467371091968
$ws.Range("F11").Value = 1.06608081694957
$ws.Range("I11").Value = 1.046755016790428
$ws.Range("J11").Value = 1.054521770931791
$ws.Range("K11").Value = 1.057436062386422
$ws.Range("L11").Value = 1.058963775056008
$ws.Range("M11").Value = 1.069539874573744
$ws.Range("N11").Value = 1.021941700788583
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04781463387169
$ws.Range("D12").Value = 1.053756096132556
$ws.Range("E12").Value = 1.055266789458243
$ws.Range("F12").Value = 1.065870634775459
$ws.Range("I12").Value = 1.046693006604673
$ws.Range("J12").Value = 1.054375984006474
$ws.Range("K12").Value = 1.057302706417481
$ws.Range("L12").Value = 1.058807942071501
$ws.Range("M12").Value = 1.069373959068503
$ws.Range("N12").Value = 1.021892232749897
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.047864215103246
$ws.Range("D13").Value = 1.053794293974439
$ws.Range("E13").Value = 1.055309807064248
$ws.Range("F13").Value = 1.065915711385696
$ws.Range("I13").Value = 1.046706314773047
$ws.Range("J13").Value = 1.054407255742546
$ws.Range("K13").Value = 1.05733131246557
$ws.Range("L13").Value = 1.058841366464805
$ws.Range("M13").Value = 1.069409545790509
$ws.Range("N13").Value = 1.021902844422718
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.048026692649276
$ws.Range("D14").Value = 1.053919471010528
$ws.Range("E14").Value = 1.055450787421692
$ws.Range("F14").Value = 1.066063439540546
$ws.Range("I14").Value = 1.046749894063316
$ws.Range("J14").Value = 1.054509720116217
$ws.Range("K14").Value = 1.057425039475687
$ws.Range("L14").Value = 1.058950892812968
$ws.Range("M14").Value = 1.069526158723431
$ws.Range("N14").Value = 1.021937612020993
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.048126816449177
$ws.Range("D15").Value = 1.053996611206965
$ws.Range("E15").Value = 1.055537672989512
$ws.Range("F15").Value = 1.066154483655837
$ws.Range("I15").Value = 1.046776724858631
$ws.Range("J15").Value = 1.054572851938441
$ws.Range("K15").Value = 1.05748278559366
$ws.Range("L15").Value = 1.059018382356402
$ws.Range("M15").Value = 1.069598015733253
$ws.Range("N15").Value = 1.021959031690004
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.048709823953316
$ws.Range("D16").Value = 1.054445820622854
$ws.Range("E16").Value = 1.056043733261761
$ws.Range("F16").Value = 1.066684763785845
$ws.Range("I16").Value = 1.046932588224504
$ws.Range("J16").Value = 1.054940307720362
$ws.Range("K16").Value = 1.057818858725484
$ws.Range("L16").Value = 1.059411302725916
$ws.Range("M16").Value = 1.070016374695413
$ws.Range("N16").Value = 1.02208367565882
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.049075742829234
$ws.Range("D17").Value = 1.054727790450346
$ws.Range("E17").Value = 1.056361477341674
$ws.Range("F17").Value = 1.067017713767773
$ws.Range("I17").Value = 1.047030089421687
$ws.Range("J17").Value = 1.055170803252086
$ws.Range("K17").Value = 1.058029637456361
$ws.Range("L17").Value = 1.059657860168731
$ws.Range("M17").Value = 1.070278904903611
$ws.Range("N17").Value = 1.022161836705737
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.04928925373523
$ws.Range("D18").Value = 1.05489232760549
$ws.Range("E18").Value = 1.056546921798131
$ws.Range("F18").Value = 1.067212032383023
$ws.Range("I18").Value = 1.047086863520616
$ws.Range("J18").Value = 1.055305247119108
$ws.Range("K18").Value = 1.058152569566924
$ws.Range("L18").Value = 1.059801704583598
$ws.Range("M18").Value = 1.070432071623879
$ws.Range("N18").Value = 1.022207417630575
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.049362068395422
$ws.Range("D19").Value = 1.054948442167421
$ws.Range("E19").Value = 1.056610172125973
$ws.Range("F19").Value = 1.067278309383777
$ws.Range("I19").Value = 1.047106205600526
$ws.Range("J19").Value = 1.055351088922487
$ws.Range("K19").Value = 1.058194484245311
$ws.Range("L19").Value = 1.059850757138962
$ws.Range("M19").Value = 1.070484303801418
$ws.Range("N19").Value = 1.022222957992137
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.049036475263333
$ws.Range("D20").Value = 1.054697530615994
$ws.Range("E20").Value = 1.056327375049668
$ws.Range("F20").Value = 1.066981979527957
$ws.Range("I20").Value = 1.047019638463641
$ws.Range("J20").Value = 1.055146073283493
$ws.Range("K20").Value = 1.058007024080606
$ws.Range("L20").Value = 1.0596314036079
$ws.Range("M20").Value = 1.07025073402982
$ws.Range("N20").Value = 1.0221534517022
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.047978843721812
$ws.Range("D21").Value = 1.05388260647132
$ws.Range("E21").Value = 1.05540926744333
$ws.Range("F21").Value = 1.066019932294156
$ws.Range("I21").Value = 1.04673706518132
$ws.Range("J21").Value = 1.054479546880881
$ws.Range("K21").Value = 1.057397439660733
$ws.Range("L21").Value = 1.058918638625819
$ws.Range("M21").Value = 1.069491817466089
$ws.Range("N21").Value = 1.021927374203143
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.047314567346249
$ws.Range("D22").Value = 1.05337086287279
$ws.Range("E22").Value = 1.054833016060431
$ws.Range("F22").Value = 1.065416098219452
$ws.Range("I22").Value = 1.046558533219699
$ws.Range("J22").Value = 1.054060480477318
$ws.Range("K22").Value = 1.05701407336862
$ws.Range("L22").Value = 1.058470788508363
$ws.Range("M22").Value = 1.069015001782791
$ws.Range("N22").Value = 1.021785151530688
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.047666646202252
$ws.Range("D23").Value = 1.053642087458144
$ws.Range("E23").Value = 1.055138402568327
$ws.Range("F23").Value = 1.065736102647467
$ws.Range("I23").Value = 1.046653258358549
$ws.Range("J23").Value = 1.054282634640241
$ws.Range("K23").Value = 1.057217311852696
$ws.Range("L23").Value = 1.058708174091633
$ws.Range("M23").Value = 1.069267737616938
$ws.Range("N23").Value = 1.021860553775557
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.04905421834949
$ws.Range("D24").Value = 1.054711203519833
$ws.Range("E24").Value = 1.056342784068709
$ws.Range("F24").Value = 1.066998125939104
$ws.Range("I24").Value = 1.047024361100576
$ws.Range("J24").Value = 1.055157247692786
$ws.Range("K24").Value = 1.058017242127562
$ws.Range("L24").Value = 1.059643358091338
$ws.Range("M24").Value = 1.070263463120014
$ws.Range("N24").Value = 1.022157240552512
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.050667419744727
$ws.Range("D25").Value = 1.055954551417772
$ws.Range("E25").Value = 1.05774469973348
$ws.Range("F25").Value = 1.068467120237040
$ws.Range("I25").Value = 1.04745120961271
$ws.Range("J25").Value = 1.056172172547296
$ws.Range("K25").Value = 1.058945059252454
$ws.Range("L25").Value = 1.060729823342228
$ws.Range("M25").Value = 1.071420405109629
$ws.Range("N25").Value = 1.02250117075421
